$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.38267251156081272
$ws.Range("A2").Value = -0.0099999993713915103
$ws.Range("A3").Value = -0.063903958957496343
$ws.Range("A4").Value = -0.011999999824901408
$ws.Range("A5").Value = -0.0059999993584289868
$ws.Range("A6").Value = -0.0059999993378987426
$ws.Range("A7").Value = -0.019999999225223775
$ws.Range("A8").Value = -0.019999999221517406
$ws.Range("A9").Value = -0.0059999993301991239
$ws.Range("A10").Value = -0.0059999993288215592
$ws.Range("A11").Value = -0.0044999993407230932
$ws.Range("A12").Value = -0.0059999993288268882
$ws.Range("A13").Value = -0.0059999993305801524
$ws.Range("A14").Value = -0.011138823104215767
$ws.Range("A15").Value = -0.0059999993339250324
$ws.Range("A16").Value = 0.045609868313091884
$ws.Range("A17").Value = -0.0059999993344170832
$ws.Range("A18").Value = -0.0089999993097231723
$ws.Range("A19").Value = -0.0089999993809799506
$ws.Range("A20").Value = -0.008999999374978529
$ws.Range("A21").Value = -0.0089999993739500184
$ws.Range("A22").Value = -0.0089999993731506578
$ws.Range("A23").Value = -0.0089999993481466589
$ws.Range("A24").Value = -0.041999999074698202
$ws.Range("A25").Value = -0.041999999068885074
$ws.Range("A26").Value = -0.0070303228089798608
$ws.Range("A27").Value = -0.0059999993318449185
$ws.Range("A28").Value = -0.0059999993178498912
$ws.Range("A29").Value = -0.011999999260346783
$ws.Range("A30").Value = -0.019999999191838036
$ws.Range("A31").Value = -0.014999999224656335
$ws.Range("A32").Value = -0.018765075719480606
$ws.Range("A33").Value = -0.0059999992947021852
